# Update the fromagerie "Photo" hyperlink-cell text so the display text
# points at the fromagerie/ image folder instead of apiculture/ (the
# underlying hyperlink target URLs are left untouched, matching the
# source diff which only rewrites the shared-string text).
$wb = $excel.ActiveWorkbook
$fromagerie = $wb.Worksheets.Item("fromagerie")

$fromagerie.Range("C2").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/fromagerie/raclette.png"
$fromagerie.Range("C3").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/fromagerie/raclette.png"
$fromagerie.Range("C4").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/fromagerie/morbier.png"
$fromagerie.Range("C5").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/fromagerie/bleu.png"
$fromagerie.Range("C6").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/fromagerie/lactique.png"
$fromagerie.Range("C7").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/fromagerie/lactique.png"

# Re-select the fromagerie sheet (it becomes the active tab again) and
# move the selection to C8, mirroring the saved UI state captured in the
# workbook/sheet XML.
$fromagerie.Activate() | Out-Null
$fromagerie.Range("C8").Select() | Out-Null
